$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "SportsName"
$ws.Range("B1").Value = "Team Name"
$ws.Range("C1").Value = "Captain Name"
$ws.Range("D1").Value = "Full Name"
$ws.Range("E1").Value = "Contact Number"
$ws.Range("F1").Value = "Email "
$ws.Range("G1").Value = "Course Name"
$ws.Range("H1").Value = "Year"
$ws.Range("I1").Value = "Category "
$ws.Range("J1").Value = "UPI ID (payment proof)"

# ---- Row 2 (ArmWrestling) ----
$ws.Range("A2").Value = "ArmWrestling"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = "entry.2076395429"
$ws.Range("E2").Value = "entry.555639742"
$ws.Range("F2").Value = "entry.548555184"
$ws.Range("G2").Value = "entry.1828627585"
$ws.Range("H2").Value = "entry.2093456979"
$ws.Range("I2").Value = "entry.139421432"
$ws.Range("J2").Value = "entry.139421432"

# ---- Row 3 (Football) ----
$ws.Range("A3").Value = "Football "
$ws.Range("B3").Value = "entry.1514101190"
$ws.Range("C3").Value = "entry.2076395429"
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = "entry.555639742"
$ws.Range("F3").Value = "entry.548555184"
$ws.Range("G3").Value = "entry.1828627585"
$ws.Range("H3").Value = "entry.2093456979"
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = "entry.139421432"

# ---- Row 4 (Cricket) ----
$ws.Range("A4").Value = "Cricket "
$ws.Range("B4").Value = "entry.1514101190"
$ws.Range("C4").Value = "entry.2076395429"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = "entry.555639742"
$ws.Range("F4").Value = "entry.548555184"
$ws.Range("G4").Value = "entry.1828627585"
$ws.Range("H4").Value = "entry.2093456979"
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = "entry.139421432"

# ---- Selection ----
$ws.Range("H11").Select() | Out-Null
